$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J2").Formula = "=E2/D2*100"
$ws.Range("J3").Formula = "=E3/D3*100"

$wb.Save()
